# Auto-generated edit script applying updated market price data
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets,
# matching the scheduled runner's recalculated values.

$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = 'ALC'; Cell = 'H28'; Value = 1316 },
    @{ Sheet = 'ALC'; Cell = 'I28'; Value = 1699 },
    @{ Sheet = 'ALC'; Cell = 'J28'; Value = 1060.6666 },
    @{ Sheet = 'ALC'; Cell = 'K28'; Value = 1699 },
    @{ Sheet = 'ALC'; Cell = 'L28'; Value = 1060.6666 },
    @{ Sheet = 'ALC'; Cell = 'M28'; Value = -1214 },
    @{ Sheet = 'ALC'; Cell = 'N28'; Value = -2030.6666 },
    @{ Sheet = 'ALC'; Cell = 'H98'; Value = 979.8 },
    @{ Sheet = 'ALC'; Cell = 'J98'; Value = 983 },
    @{ Sheet = 'ALC'; Cell = 'L98'; Value = 983 },
    @{ Sheet = 'ALC'; Cell = 'N98'; Value = -3979 },
    @{ Sheet = 'ALC'; Cell = 'H113'; Value = 7893.7856 },
    @{ Sheet = 'ALC'; Cell = 'I113'; Value = 7958.909 },
    @{ Sheet = 'ALC'; Cell = 'J113'; Value = 7655 },
    @{ Sheet = 'ALC'; Cell = 'K113'; Value = 7958.909 },
    @{ Sheet = 'ALC'; Cell = 'L113'; Value = 7655 },
    @{ Sheet = 'ALC'; Cell = 'M113'; Value = -4704.909 },
    @{ Sheet = 'ALC'; Cell = 'N113'; Value = -14163 },
    @{ Sheet = 'ALC'; Cell = 'H122'; Value = 979.8 },
    @{ Sheet = 'ALC'; Cell = 'J122'; Value = 983 },
    @{ Sheet = 'ALC'; Cell = 'L122'; Value = 2949 },
    @{ Sheet = 'ALC'; Cell = 'N122'; Value = -7849 },
    @{ Sheet = 'ALC'; Cell = 'H137'; Value = 1141.7693 },
    @{ Sheet = 'ALC'; Cell = 'I137'; Value = 1058.875 },
    @{ Sheet = 'ALC'; Cell = 'J137'; Value = 1274.4 },
    @{ Sheet = 'ALC'; Cell = 'K137'; Value = 3176.625 },
    @{ Sheet = 'ALC'; Cell = 'L137'; Value = 3823.2 },
    @{ Sheet = 'ALC'; Cell = 'M137'; Value = -626.625 },
    @{ Sheet = 'ALC'; Cell = 'N137'; Value = -8923.200000000001 },
    @{ Sheet = 'ARM'; Cell = 'H28'; Value = 10195 },
    @{ Sheet = 'ARM'; Cell = 'I28'; Value = 10195 },
    @{ Sheet = 'ARM'; Cell = 'K28'; Value = 10195 },
    @{ Sheet = 'ARM'; Cell = 'M28'; Value = -10003 },
    @{ Sheet = 'ARM'; Cell = 'H32'; Value = 3992.0417 },
    @{ Sheet = 'ARM'; Cell = 'I32'; Value = 3740.652 },
    @{ Sheet = 'ARM'; Cell = 'K32'; Value = 3740.652 },
    @{ Sheet = 'ARM'; Cell = 'M32'; Value = -3453.652 },
    @{ Sheet = 'ARM'; Cell = 'H74'; Value = 999 },
    @{ Sheet = 'ARM'; Cell = 'I74'; Value = 999 },
    @{ Sheet = 'ARM'; Cell = 'K74'; Value = 999 },
    @{ Sheet = 'ARM'; Cell = 'M74'; Value = -125 },
    @{ Sheet = 'ARM'; Cell = 'H77'; Value = 999 },
    @{ Sheet = 'ARM'; Cell = 'I77'; Value = 999 },
    @{ Sheet = 'ARM'; Cell = 'K77'; Value = 4995 },
    @{ Sheet = 'ARM'; Cell = 'M77'; Value = -627 },
    @{ Sheet = 'ARM'; Cell = 'H99'; Value = 10195 },
    @{ Sheet = 'ARM'; Cell = 'I99'; Value = 10195 },
    @{ Sheet = 'ARM'; Cell = 'K99'; Value = 10195 },
    @{ Sheet = 'ARM'; Cell = 'M99'; Value = -7200 },
    @{ Sheet = 'ARM'; Cell = 'H102'; Value = 3848 },
    @{ Sheet = 'ARM'; Cell = 'I102'; Value = 3848 },
    @{ Sheet = 'ARM'; Cell = 'K102'; Value = 3848 },
    @{ Sheet = 'ARM'; Cell = 'M102'; Value = -2226 },
    @{ Sheet = 'ARM'; Cell = 'H110'; Value = 2426.8572 },
    @{ Sheet = 'ARM'; Cell = 'I110'; Value = 2426.8572 },
    @{ Sheet = 'ARM'; Cell = 'K110'; Value = 2426.8572 },
    @{ Sheet = 'ARM'; Cell = 'M110'; Value = -381.8571999999999 },
    @{ Sheet = 'ARM'; Cell = 'H132'; Value = 822.2857 },
    @{ Sheet = 'ARM'; Cell = 'I132'; Value = 822.2857 },
    @{ Sheet = 'ARM'; Cell = 'K132'; Value = 2466.8571 },
    @{ Sheet = 'ARM'; Cell = 'M132'; Value = 63.14289999999983 },
    @{ Sheet = 'BSM'; Cell = 'H86'; Value = 4950.25 },
    @{ Sheet = 'BSM'; Cell = 'J86'; Value = 5249.3335 },
    @{ Sheet = 'BSM'; Cell = 'L86'; Value = 5249.3335 },
    @{ Sheet = 'BSM'; Cell = 'N86'; Value = -7495.3335 },
    @{ Sheet = 'BSM'; Cell = 'H89'; Value = 4950.25 },
    @{ Sheet = 'BSM'; Cell = 'J89'; Value = 5249.3335 },
    @{ Sheet = 'BSM'; Cell = 'L89'; Value = 26246.6675 },
    @{ Sheet = 'BSM'; Cell = 'N89'; Value = -37478.6675 },
    @{ Sheet = 'BSM'; Cell = 'H107'; Value = 600 },
    @{ Sheet = 'BSM'; Cell = 'I107'; Value = 600 },
    @{ Sheet = 'BSM'; Cell = 'K107'; Value = 600 },
    @{ Sheet = 'BSM'; Cell = 'M107'; Value = 1320 },
    @{ Sheet = 'BSM'; Cell = 'H134'; Value = 1755.1111 },
    @{ Sheet = 'BSM'; Cell = 'I134'; Value = 1755.1111 },
    @{ Sheet = 'BSM'; Cell = 'J134'; Value = 0 },
    @{ Sheet = 'BSM'; Cell = 'K134'; Value = 5265.3333 },
    @{ Sheet = 'BSM'; Cell = 'L134'; Value = 0 },
    @{ Sheet = 'BSM'; Cell = 'M134'; Action = 'Clear' },
    @{ Sheet = 'BSM'; Cell = 'N134'; Value = -2730.3333 },
    @{ Sheet = 'CRP'; Cell = 'H99'; Value = 6405.778 },
    @{ Sheet = 'CRP'; Cell = 'I99'; Value = 6543.143 },
    @{ Sheet = 'CRP'; Cell = 'K99'; Value = 6543.143 },
    @{ Sheet = 'CRP'; Cell = 'M99'; Value = -5045.143 },
    @{ Sheet = 'CRP'; Cell = 'H105'; Value = 5953.5264 },
    @{ Sheet = 'CRP'; Cell = 'I105'; Value = 6173.1665 },
    @{ Sheet = 'CRP'; Cell = 'J105'; Value = 2000 },
    @{ Sheet = 'CRP'; Cell = 'K105'; Value = 6173.1665 },
    @{ Sheet = 'CRP'; Cell = 'L105'; Value = 2000 },
    @{ Sheet = 'CRP'; Cell = 'M105'; Value = -4426.1665 },
    @{ Sheet = 'CRP'; Cell = 'N105'; Value = -5494 },
    @{ Sheet = 'CRP'; Cell = 'H107'; Value = 343.9524 },
    @{ Sheet = 'CRP'; Cell = 'I107'; Value = 315.70587 },
    @{ Sheet = 'CRP'; Cell = 'J107'; Value = 464 },
    @{ Sheet = 'CRP'; Cell = 'K107'; Value = 315.70587 },
    @{ Sheet = 'CRP'; Cell = 'L107'; Value = 464 },
    @{ Sheet = 'CRP'; Cell = 'M107'; Value = 1604.29413 },
    @{ Sheet = 'CRP'; Cell = 'N107'; Value = -4304 },
    @{ Sheet = 'CRP'; Cell = 'H122'; Value = 1971.9166 },
    @{ Sheet = 'CRP'; Cell = 'I122'; Value = 1787.5454 },
    @{ Sheet = 'CRP'; Cell = 'K122'; Value = 5362.6362 },
    @{ Sheet = 'CRP'; Cell = 'M122'; Value = -2912.6362 },
    @{ Sheet = 'CRP'; Cell = 'H126'; Value = 6405.778 },
    @{ Sheet = 'CRP'; Cell = 'I126'; Value = 6543.143 },
    @{ Sheet = 'CRP'; Cell = 'K126'; Value = 19629.429 },
    @{ Sheet = 'CRP'; Cell = 'M126'; Value = -17159.429 },
    @{ Sheet = 'CRP'; Cell = 'H132'; Value = 2477.3333 },
    @{ Sheet = 'CRP'; Cell = 'I132'; Value = 2328.9412 },
    @{ Sheet = 'CRP'; Cell = 'J132'; Value = 5000 },
    @{ Sheet = 'CRP'; Cell = 'K132'; Value = 6986.823600000001 },
    @{ Sheet = 'CRP'; Cell = 'L132'; Value = 15000 },
    @{ Sheet = 'CRP'; Cell = 'M132'; Value = -4456.823600000001 },
    @{ Sheet = 'CRP'; Cell = 'N132'; Value = -20060 },
    @{ Sheet = 'CUL'; Cell = 'H12'; Value = 419.25 },
    @{ Sheet = 'CUL'; Cell = 'I12'; Value = 0 },
    @{ Sheet = 'CUL'; Cell = 'J12'; Value = 419.25 },
    @{ Sheet = 'CUL'; Cell = 'K12'; Value = 0 },
    @{ Sheet = 'CUL'; Cell = 'L12'; Action = 'Clear' },
    @{ Sheet = 'CUL'; Cell = 'M12'; Value = 1257.75 },
    @{ Sheet = 'CUL'; Cell = 'N12'; Value = -1603.75 },
    @{ Sheet = 'CUL'; Cell = 'H68'; Value = 1503 },
    @{ Sheet = 'CUL'; Cell = 'J68'; Value = 1503 },
    @{ Sheet = 'CUL'; Cell = 'L68'; Value = 4509 },
    @{ Sheet = 'CUL'; Cell = 'N68'; Value = -6131 },
    @{ Sheet = 'CUL'; Cell = 'H71'; Value = 1503 },
    @{ Sheet = 'CUL'; Cell = 'J71'; Value = 1503 },
    @{ Sheet = 'CUL'; Cell = 'L71'; Value = 13527 },
    @{ Sheet = 'CUL'; Cell = 'N71'; Value = -21639 },
    @{ Sheet = 'CUL'; Cell = 'H131'; Value = 891.36365 },
    @{ Sheet = 'CUL'; Cell = 'J131'; Value = 888.53845 },
    @{ Sheet = 'CUL'; Cell = 'L131'; Value = 2665.61535 },
    @{ Sheet = 'CUL'; Cell = 'N131'; Value = -12745.61535 },
    @{ Sheet = 'GSM'; Cell = 'H102'; Value = 2459.6 },
    @{ Sheet = 'GSM'; Cell = 'I102'; Value = 2272.2856 },
    @{ Sheet = 'GSM'; Cell = 'K102'; Value = 2272.2856 },
    @{ Sheet = 'GSM'; Cell = 'M102'; Value = -650.2856000000002 },
    @{ Sheet = 'GSM'; Cell = 'H122'; Value = 2074.1 },
    @{ Sheet = 'GSM'; Cell = 'I122'; Value = 1891.5 },
    @{ Sheet = 'GSM'; Cell = 'J122'; Value = 2348 },
    @{ Sheet = 'GSM'; Cell = 'K122'; Value = 5674.5 },
    @{ Sheet = 'GSM'; Cell = 'L122'; Value = 7044 },
    @{ Sheet = 'GSM'; Cell = 'M122'; Value = -3224.5 },
    @{ Sheet = 'GSM'; Cell = 'N122'; Value = -11944 },
    @{ Sheet = 'GSM'; Cell = 'H132'; Value = 3271.5 },
    @{ Sheet = 'GSM'; Cell = 'I132'; Value = 3271.5 },
    @{ Sheet = 'GSM'; Cell = 'K132'; Value = 9814.5 },
    @{ Sheet = 'GSM'; Cell = 'M132'; Value = -7284.5 },
    @{ Sheet = 'LTW'; Cell = 'H40'; Value = 2994 },
    @{ Sheet = 'LTW'; Cell = 'I40'; Value = 2994 },
    @{ Sheet = 'LTW'; Cell = 'K40'; Value = 2994 },
    @{ Sheet = 'LTW'; Cell = 'M40'; Value = -2858 },
    @{ Sheet = 'LTW'; Cell = 'H110'; Value = 48782 },
    @{ Sheet = 'LTW'; Cell = 'J110'; Value = 48782 },
    @{ Sheet = 'LTW'; Cell = 'L110'; Value = 48782 },
    @{ Sheet = 'LTW'; Cell = 'N110'; Value = -56962 },
    @{ Sheet = 'WVR'; Cell = 'H107'; Value = 997.5 },
    @{ Sheet = 'WVR'; Cell = 'I107'; Value = 997.5 },
    @{ Sheet = 'WVR'; Cell = 'J107'; Value = 0 },
    @{ Sheet = 'WVR'; Cell = 'K107'; Value = 2992.5 },
    @{ Sheet = 'WVR'; Cell = 'L107'; Value = 0 },
    @{ Sheet = 'WVR'; Cell = 'M107'; Action = 'Clear' },
    @{ Sheet = 'WVR'; Cell = 'N107'; Value = -1072.5 },
    @{ Sheet = 'WVR'; Cell = 'H122'; Value = 1262.25 },
    @{ Sheet = 'WVR'; Cell = 'I122'; Value = 1262.25 },
    @{ Sheet = 'WVR'; Cell = 'J122'; Value = 0 },
    @{ Sheet = 'WVR'; Cell = 'K122'; Value = 3786.75 },
    @{ Sheet = 'WVR'; Cell = 'L122'; Value = 0 },
    @{ Sheet = 'WVR'; Cell = 'M122'; Action = 'Clear' },
    @{ Sheet = 'WVR'; Cell = 'N122'; Value = -1336.75 },
    @{ Sheet = 'WVR'; Cell = 'H126'; Value = 3592.125 },
    @{ Sheet = 'WVR'; Cell = 'I126'; Value = 2308.1667 },
    @{ Sheet = 'WVR'; Cell = 'K126'; Value = 6924.500100000001 },
    @{ Sheet = 'WVR'; Cell = 'M126'; Value = -4454.500100000001 },
    @{ Sheet = 'WVR'; Cell = 'H132'; Value = 2994.375 },
    @{ Sheet = 'WVR'; Cell = 'I132'; Value = 2742.5 },
    @{ Sheet = 'WVR'; Cell = 'K132'; Value = 8227.5 },
    @{ Sheet = 'WVR'; Cell = 'M132'; Value = -5697.5 },
)

$sheetCache = @{}
foreach ($change in $changes) {
    if (-not $sheetCache.ContainsKey($change.Sheet)) {
        $sheetCache[$change.Sheet] = $wb.Worksheets.Item($change.Sheet)
    }
    $ws = $sheetCache[$change.Sheet]
    $rng = $ws.Range($change.Cell)
    if ($change.Action -eq 'Clear') {
        $rng.ClearContents()
    } else {
        $rng.Value = $change.Value
    }
}

Write-Host "Applied" $changes.Count "cell updates across" $sheetCache.Count "sheets"
